# Add "Wins" / "Losses" / "Ties" season-record columns (AC:AE) to Sheet1.
# The existing sheet uses columns A:AB (header row 1, data rows 2:43).
# New columns get the same header style as the existing headers (copied
# from AB1) and each data row gets the team's season record: 76 wins,
# 86 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 43

# --- Headers (row 1) ------------------------------------------------------
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows (2 .. lastRow) ---------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 76   # AC - Wins
    $ws.Cells.Item($r, 30).Value = 86   # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE - Ties
}
